# Apply cell value updates from the commit diff.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the source file, where these "numbers" are inline strings,
# e.g. multi-dot prices like "57.426.24" or zero-padded values like "2.00").
# Resetting .Style to "Normal" afterwards keeps the cell on the default
# style index (0) instead of the auto-generated "quote prefix" style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $range = $ws.Range($ref)
    $range.Value = "'" + $val
    $range.Style = 'Normal'
}

Set-TextValue "D2" "57.314.18"
Set-TextValue "E2" "  -7.27%  "
Set-TextValue "D3" "2.878.38"
Set-TextValue "E3" "  -5.81%  "
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "550.96"
Set-TextValue "E5" "  -6.11%  "
Set-TextValue "D6" "121.11"
Set-TextValue "E6" "  -7.11%  "
Set-TextValue "E7" "  +0.25%  "
Set-TextValue "D8" "2.870.84"
Set-TextValue "E8" "  -5.97%  "
Set-TextValue "E9" "  -2.85%  "
Set-TextValue "E10" "  -11.11%  "
Set-TextValue "D11" "4.74"
Set-TextValue "E11" "  -10.19%  "
Set-TextValue "D12" "0.431"
Set-TextValue "E12" "  -2.46%  "
Set-TextValue "E13" "  -11.86%  "
Set-TextValue "D14" "31.24"
Set-TextValue "E14" "  -7.50%  "
Set-TextValue "E15" "  -0.87%  "
Set-TextValue "D16" "3.359.58"
Set-TextValue "E16" "  -5.56%  "
Set-TextValue "D17" "2.878.78"
Set-TextValue "E17" "  -5.56%  "
Set-TextValue "D18" "57.295.87"
Set-TextValue "E18" "  -7.44%  "
Set-TextValue "D19" "6.41"
Set-TextValue "E19" "  +0.13%  "
Set-TextValue "D20" "407.25"
Set-TextValue "E20" "  -9.59%  "
Set-TextValue "D21" "12.72"
Set-TextValue "E21" "  -6.24%  "
Set-TextValue "E22" "  -3.83%  "
Set-TextValue "D23" "6.70"
Set-TextValue "E23" "  -9.18%  "
Set-TextValue "D24" "12.51"
Set-TextValue "E24" "  -3.02%  "
Set-TextValue "D25" "76.53"
Set-TextValue "E25" "  -5.92%  "
Set-TextValue "D26" "0.997"
Set-TextValue "E26" "  -0.34%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.03%  "
Set-TextValue "E28" "  -4.77%  "
Set-TextValue "E29" "  -4.58%  "
Set-TextValue "D30" "1.89"
Set-TextValue "E30" "  -6.36%  "
Set-TextValue "D31" "6.02"
Set-TextValue "E31" "  -6.88%  "
Set-TextValue "D32" "24.51"
Set-TextValue "E32" "  -5.63%  "
Set-TextValue "D33" "0.0949"
Set-TextValue "E33" "  -3.09%  "
Set-TextValue "B34" "Stacks"
Set-TextValue "C34" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D34" "2.00"
Set-TextValue "E34" "  -14.33%  "
Set-TextValue "B35" "Filecoin"
Set-TextValue "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D35" "5.34"
Set-TextValue "E35" "  -6.91%  "
Set-TextValue "B36" "Mantle"
Set-TextValue "C36" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D36" "0.893"
Set-TextValue "E36" "  -8.71%  "
Set-TextValue "D37" "48.35"
Set-TextValue "E37" "  -4.14%  "
Set-TextValue "D38" "8.35"
Set-TextValue "E38" "  +4.82%  "
Set-TextValue "D39" "0.0₃0608"
Set-TextValue "E39" "  -12.64%  "
Set-TextValue "E40" "  -9.26%  "
Set-TextValue "E41" "  -3.82%  "
Set-TextValue "D42" "2.591.57"
Set-TextValue "E42" "  -3.94%  "
Set-TextValue "D43" "355.96"
Set-TextValue "E43" "  -6.97%  "
Set-TextValue "E44" "  +0.01%  "
Set-TextValue "E45" "  -7.79%  "
Set-TextValue "D46" "117.20"
Set-TextValue "E46" "  -5.60%  "
Set-TextValue "E47" "  -5.79%  "
Set-TextValue "E48" "  -3.28%  "
Set-TextValue "E49" "  -5.48%  "
Set-TextValue "D50" "22.34"
Set-TextValue "E50" "  -7.28%  "
Set-TextValue "E51" "  -7.75%  "
